$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2024-05-15", "12:52:01", "-", "-", "-", "-", "Error en sensor de salida", "12:52:05"),
    @("2024-05-15", "12:52:09", "-", "-", "-", "-", "Sensor de PCB detecta que hay placa cuando no la hay", "12:52:13"),
    @("2024-05-15", "12:52:16", "-", "-", "-", "-", "Fallo dispensaci" + [char]0x00F3 + "n glue", "12:52:23"),
    @("2024-05-15", "12:52:35", "-", "-", "-", "-", "Soldadura defectuosa", "12:52:39")
)

$startRow = 66
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le 8; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $text = $rowData[$c - 1]
        if ($c -eq 1) {
            # Column A holds a date-shaped string ("2024-05-15"); without a
            # text format Excel silently reinterprets it as a date serial.
            # Force text, write the literal, then drop the format again so
            # the cell ends up with the default (unstyled) text cell that
            # the rest of the sheet uses.
            $cell.NumberFormat = "@"
            $cell.Value = $text
            $cell.ClearFormats()
        } else {
            $cell.Value = $text
        }
    }
}
